$d = $word.ActiveDocument

$pairs = @(
    ,@('2023-03-21 Tuesday', '2023-03-22 Wednesday')
    ,@('3+82=', '54-21=')
    ,@('83-74=', '92-54=')
    ,@('50-12=', '75-48=')
    ,@('83-50=', '84-82=')
    ,@('27+70=', '57-24=')
    ,@('38+61=', '53+41=')
    ,@('63+5=', '39-25=')
    ,@('72+6=', '37+62=')
    ,@('80-33=', '73-38=')
    ,@('93-11=', '4+38=')
    ,@('17+15=', '51-19=')
    ,@('44-40=', '99-40=')
    ,@('30+32=', '86-24=')
    ,@('23+20=', '62-30=')
    ,@('84-21=', '69+29=')
    ,@('36+45=', '75-71=')
    ,@('9+56=', '16+66=')
    ,@('0+89=', '70+7=')
    ,@('57+11=', '41+15=')
    ,@('94-9=', '95-77=')
    ,@('48-14=', '1+2=')
    ,@('58-46=', '97-74=')
    ,@('2+24=', '43-32=')
    ,@('17+6=', '49+12=')
    ,@('43+41=', '13+8=')
    ,@('48+38=', '0+1=')
    ,@('25+36=', '19-17=')
    ,@('96-32=', '97-70=')
    ,@('87-42=', '91-81=')
    ,@('26+24=', '94-43=')
    ,@('34-26=', '95-66=')
    ,@('49-35=', '98-40=')
    ,@('90-12=', '66+4=')
    ,@('90-84=', '29+21=')
    ,@('84+12=', '7+34=')
    ,@('66-64=', '13+86=')
    ,@('85-72=', '68-27=')
    ,@('3+54=', '21-5=')
    ,@('51-12=', '47-18=')
    ,@('42+49=', '53-15=')
    ,@('16+40=', '74-10=')
    ,@('2+38=', '80+13=')
    ,@('32+40=', '52+8=')
    ,@('50+4=', '5+88=')
    ,@('79+17=', '7-2=')
    ,@('93-84=', '53+27=')
    ,@('54-48=', '36-30=')
    ,@('11+80=', '16+9=')
    ,@('39-9=', '63+16=')
    ,@('39+10=', '15+0=')
    ,@('3+75=', '97-40=')
    ,@('52+15=', '27-22=')
    ,@('14+13=', '71-59=')
    ,@('12+70=', '90-21=')
    ,@('69-59=', '70-7=')
    ,@('74-4=', '48+34=')
    ,@('71-43=', '76-41=')
    ,@('69-39=', '34+54=')
    ,@('46-7=', '84-15=')
    ,@('3+33=', '2+41=')
    ,@('68-11=', '69-36=')
    ,@('59-32=', '81-12=')
    ,@('78-46=', '22+76=')
    ,@('4+66=', '86-4=')
    ,@('17+62=', '73-6=')
    ,@('26+25=', '35-7=')
    ,@('50-17=', '77-70=')
    ,@('61+30=', '20+44=')
    ,@('68+5=', '90-57=')
    ,@('83-5=', '92-56=')
    ,@('87-62=', '25+67=')
    ,@('80-50=', '96-6=')
    ,@('30+25=', '59+14=')
    ,@('86-26=', '76-75=')
    ,@('88+0=', '94-29=')
    ,@('68-40=', '90-18=')
    ,@('99-88=', '79-42=')
    ,@('25+29=', '60-45=')
    ,@('86+7=', '15+31=')
    ,@('6+65=', '58+41=')
    ,@('63+29=', '39+25=')
    ,@('93-65=', '82-30=')
    ,@('5+70=', '11+21=')
    ,@('31-30=', '65-21=')
    ,@('91+0=', '12+56=')
    ,@('79-57=', '30+7=')
    ,@('21+50=', '45+8=')
    ,@('51-23=', '73-50=')
    ,@('99-77=', '46-11=')
    ,@('14+57=', '99-60=')
    ,@('18-18=', '21+41=')
    ,@('33-6=', '18+39=')
    ,@('36-25=', '71-24=')
    ,@('25+60=', '33+6=')
    ,@('33+4=', '62-27=')
    ,@('58+13=', '2+64=')
    ,@('10+55=', '42-0=')
    ,@('79-20=', '58-47=')
    ,@('18+16=', '44-30=')
    ,@('6+17=', '75-46=')
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "done"
